# Update crypto price/volume figures to the latest scrape (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values that are purely numeric-looking (e.g. "214.91") need a leading
# apostrophe so Excel stores them as text (matching the workbook's existing text cells)
# instead of silently converting them to numbers.

$ws.Range("D2").Value = '26.685.90'
$ws.Range("E2").Value = '  +0.79%  '
$ws.Range("D3").Value = '1.644.06'
$ws.Range("E3").Value = '  +1.06%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '''214.91'
$ws.Range("E5").Value = '  +0.90%  '
$ws.Range("E6").Value = '  +1.86%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("E8").Value = '  +0.97%  '
$ws.Range("D9").Value = '''0.0626'
$ws.Range("E9").Value = '  +0.91%  '
$ws.Range("D10").Value = '''19.07'
$ws.Range("E10").Value = '  +0.51%  '
$ws.Range("D11").Value = '''0.0845'
$ws.Range("E11").Value = '  +0.60%  '
$ws.Range("D12").Value = '1.872.42'
$ws.Range("E12").Value = '  +1.00%  '
$ws.Range("D13").Value = '1.668.65'
$ws.Range("E13").Value = '  +2.01%  '
$ws.Range("E14").Value = '  +1.69%  '
$ws.Range("E15").Value = '  +1.65%  '
$ws.Range("D16").Value = '''65.03'
$ws.Range("E16").Value = '  +1.81%  '
$ws.Range("D17").Value = '26.701.51'
$ws.Range("E17").Value = '  +0.79%  '
$ws.Range("D18").Value = '0.0₃0743'
$ws.Range("E18").Value = '  +0.51%  '
$ws.Range("D19").Value = '''215.86'
$ws.Range("E19").Value = '  +0.52%  '
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("E21").Value = '  +1.06%  '
$ws.Range("E22").Value = '  +0.97%  '
$ws.Range("E23").Value = '  +1.98%  '
$ws.Range("E24").Value = '  +13.75%  '
$ws.Range("D25").Value = '''145.31'
$ws.Range("E25").Value = '  -2.41%  '
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("E28").Value = '  +4.81%  '
$ws.Range("E30").Value = '  +1.36%  '
$ws.Range("D31").Value = '''1.17'
$ws.Range("E31").Value = '  +1.35%  '
$ws.Range("E32").Value = '  +1.59%  '
$ws.Range("D33").Value = '''3.00'
$ws.Range("E33").Value = '  +2.31%  '
$ws.Range("D34").Value = '1.282.76'
$ws.Range("E34").Value = '  +5.14%  '
$ws.Range("D35").Value = '''1.53'
$ws.Range("E35").Value = '  +2.62%  '
$ws.Range("E36").Value = '  +1.01%  '
$ws.Range("E37").Value = '  +2.72%  '
$ws.Range("E38").Value = '  +6.21%  '
$ws.Range("E39").Value = '  +4.10%  '
$ws.Range("D40").Value = '''1.00'
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("E41").Value = '  +2.18%  '
$ws.Range("D43").Value = '''5.43'
$ws.Range("E43").Value = '  +1.32%  '
$ws.Range("D44").Value = '1.782.33'
$ws.Range("E44").Value = '  +1.03%  '
$ws.Range("D45").Value = '''91.74'
$ws.Range("E45").Value = '  -1.28%  '
$ws.Range("D46").Value = '''59.27'
$ws.Range("E46").Value = '  +8.26%  '
$ws.Range("E47").Value = '  +1.38%  '
$ws.Range("D48").Value = '''0.0515'
$ws.Range("E48").Value = '  +0.88%  '
$ws.Range("D49").Value = '''7.76'
$ws.Range("E49").Value = '  +1.55%  '
$ws.Range("D50").Value = '''0.0965'
$ws.Range("E50").Value = '  +2.16%  '
$ws.Range("D51").Value = '''0.406'
$ws.Range("E51").Value = '  -0.51%  '
